# Femacal de La Calera - Pepino ensalada: add a new weekly price record.
# The new record is inserted as row 253 (pushing every existing record
# at/after the old row 253 down by one row), growing the used range
# from A1:R275 to A1:R276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 253; Excel shifts rows 253:275 down
# to 254:276 and extends the sheet dimension automatically.
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A253").Value = 3
$ws.Range("B253").Value = "Femacal de La Calera"
$ws.Range("C253").Value = "Coquimbo"
$ws.Range("D253").Value = 44578
$ws.Range("E253").Value = 5
$ws.Range("F253").Value = 100112043
$ws.Range("G253").Value = "Pepino ensalada"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 105
$ws.Range("K253").Value = 11000
$ws.Range("L253").Value = 12000
$ws.Range("M253").Value = 11476
$ws.Range("N253").Value = "$/caja 70 unidades"
$ws.Range("O253").Value = "Limache"
$ws.Range("P253").Value = 164
$ws.Range("Q253").Value = 70
$ws.Range("R253").Value = "Hortaliza"
